$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 6 (the spurious divider row "grandes regiões e unidades
# da federação" that had a label but no data) and shift everything below it
# up by one row.
$ws.Rows.Item(6).Delete()
